$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.919643
$ws.Range("H2").Value = 53.75892899999999
$ws.Range("I2").Value = 0.8982899767221961
$ws.Range("J2").Value = 0.8982899767221962
$ws.Range("M2").Value = 0.05661333333333334
$ws.Range("N2").Value = 0.16984
$ws.Range("O2").Value = 0.0204119846136133
$ws.Range("P2").Value = 0.02041198461361329
$ws.Range("Q2").Value = 1.014490722373333
$ws.Range("R2").Value = 9.130416501360001
$ws.Range("S2").Value = 0.01833588118341651
$ws.Range("T2").Value = 0.01833588118341651

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.919643
$ws.Range("H3").Value = 53.75892899999999
$ws.Range("I3").Value = 0.8982899767221961
$ws.Range("J3").Value = 0.8982899767221962
$ws.Range("N3").Value = 0.8341160000000001
$ws.Range("O3").Value = 0.100247073468963
$ws.Range("P3").Value = 0.1002470734689629
$ws.Range("Q3").Value = 4.982353646862666
$ws.Range("R3").Value = 44.841182821764
$ws.Range("S3").Value = 0.09005094129290302
$ws.Range("T3").Value = 0.09005094129290302

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.919643
$ws.Range("H4").Value = 53.75892899999999
$ws.Range("I4").Value = 0.8982899767221961
$ws.Range("J4").Value = 0.8982899767221962
$ws.Range("M4").Value = 2.438882
$ws.Range("N4").Value = 7.316646
$ws.Range("O4").Value = 0.8793409419174237
$ws.Range("P4").Value = 0.8793409419174237
$ws.Range("Q4").Value = 43.70389475912599
$ws.Range("R4").Value = 393.335052832134
$ws.Range("S4").Value = 0.7899031542458765
$ws.Range("T4").Value = 0.7899031542458766

$ws.Range("I5").Value = 0.06812533974785755
$ws.Range("J5").Value = 0.06812533974785755
$ws.Range("M5").Value = 0.05661333333333334
$ws.Range("N5").Value = 0.16984
$ws.Range("O5").Value = 0.0204119846136133
$ws.Range("P5").Value = 0.02041198461361329
$ws.Range("Q5").Value = 0.07693787855111112
$ws.Range("R5").Value = 0.69244090696
$ws.Range("S5").Value = 0.001390573386730447
$ws.Range("T5").Value = 0.001390573386730446

$ws.Range("I6").Value = 0.06812533974785755
$ws.Range("J6").Value = 0.06812533974785755
$ws.Range("N6").Value = 0.8341160000000001
$ws.Range("O6").Value = 0.100247073468963
$ws.Range("P6").Value = 0.1002470734689629
$ws.Range("Q6").Value = 0.3778563089115556
$ws.Range("S6").Value = 0.006829365938801538
$ws.Range("T6").Value = 0.006829365938801537

$ws.Range("I7").Value = 0.06812533974785755
$ws.Range("J7").Value = 0.06812533974785755
$ws.Range("M7").Value = 2.438882
$ws.Range("N7").Value = 7.316646
$ws.Range("O7").Value = 0.8793409419174237
$ws.Range("P7").Value = 0.8793409419174237
$ws.Range("Q7").Value = 3.314456084252666
$ws.Range("R7").Value = 29.830104758274
$ws.Range("S7").Value = 0.05990540042232556
$ws.Range("T7").Value = 0.05990540042232556

$ws.Range("G8").Value = 0.669968
$ws.Range("H8").Value = 2.009904
$ws.Range("I8").Value = 0.03358468352994624
$ws.Range("J8").Value = 0.03358468352994624
$ws.Range("M8").Value = 0.05661333333333334
$ws.Range("N8").Value = 0.16984
$ws.Range("O8").Value = 0.0204119846136133
$ws.Range("P8").Value = 0.02041198461361329
$ws.Range("Q8").Value = 0.03792912170666667
$ws.Range("R8").Value = 0.3413620953600001
$ws.Range("S8").Value = 0.0006855300434663344
$ws.Range("T8").Value = 0.0006855300434663345

$ws.Range("G9").Value = 0.669968
$ws.Range("H9").Value = 2.009904
$ws.Range("I9").Value = 0.03358468352994624
$ws.Range("J9").Value = 0.03358468352994624
$ws.Range("N9").Value = 0.8341160000000001
$ws.Range("O9").Value = 0.100247073468963
$ws.Range("P9").Value = 0.1002470734689629
$ws.Range("Q9").Value = 0.1862770094293334
$ws.Range("R9").Value = 1.676493084864
$ws.Range("S9").Value = 0.003366766237258391
$ws.Range("T9").Value = 0.003366766237258391

$ws.Range("G10").Value = 0.669968
$ws.Range("H10").Value = 2.009904
$ws.Range("I10").Value = 0.03358468352994624
$ws.Range("J10").Value = 0.03358468352994624
$ws.Range("M10").Value = 2.438882
$ws.Range("N10").Value = 7.316646
$ws.Range("O10").Value = 0.8793409419174237
$ws.Range("P10").Value = 0.8793409419174237
$ws.Range("Q10").Value = 1.633972895776
$ws.Range("R10").Value = 14.705756061984
$ws.Range("S10").Value = 0.02953238724922151
$ws.Range("T10").Value = 0.02953238724922152
